$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..41 (row 1 is the header: index, place, type, country, price, rent)
for ($row = 2; $row -le 41; $row++) {
    $type = $ws.Cells.Item($row, 3).Value2          # column C = type
    $price = $ws.Cells.Item($row, 5).Value2         # column E = price

    # Only "city" rows keep their real country; everything else (stay_place,
    # service_centers, ...) becomes an "Upgrade" option with no fixed country.
    if ($type -ne "city") {
        $ws.Cells.Item($row, 4).Value2 = "None"     # column D = country
    }

    # Replace the rent formula with its computed static value (30% of price)
    $ws.Cells.Item($row, 6).Value2 = (30 / 100) * $price   # column F = rent
}
